# Apply the "JE tblversion base terminer" edit:
#  - Renumber the IDJeu (H column) group ids for the existing Halo4 rows
#    so the Multi-local / Multi-online groups are distinguished from Solo.
#  - Append the Mario Kart 8, World of Warcraft and NHL 14/15 version rows
#    (rows 16-36) to the tblVersion sheet, including the new note in J20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tblVersion")

# --- renumber existing IDJeu values (rows 8-14) ---
$ws.Range("H8").Value  = 3
$ws.Range("H9").Value  = 3
$ws.Range("H10").Value = 4
$ws.Range("H12").Value = 4
$ws.Range("H13").Value = 5
$ws.Range("H14").Value = 5

# --- add the new version rows (16-36) ---
$ws.Range("A16").Value = 'MKart8SS01'
$ws.Range("B16").Value = 'Mario Kart 8 solo 1.0'
$ws.Range("C16").Value = '1er version du solo de Mario Kart 8, bug nombreux'
$ws.Range("D16").Value = 'Beta'
$ws.Range("E16").Value = 20130105
$ws.Range("F16").Value = 20140530
$ws.Range("H16").Value = 6
$ws.Range("A17").Value = 'MKart8SS02'
$ws.Range("B17").Value = 'Mario Kart 8 solo 1.1'
$ws.Range("C17").Value = '2er version du solo de Mario Kart 8, ajout mode course, mini-jeux, etc'
$ws.Range("D17").Value = 'Beta'
$ws.Range("E17").Value = 20140110
$ws.Range("F17").Value = 20140530
$ws.Range("H17").Value = 6
$ws.Range("A18").Value = 'MKart8ML01'
$ws.Range("B18").Value = 'Mario Kart 8 Multi Local 1.0'
$ws.Range("C18").Value = '1er version du Multi Local de Mario Kart 8, bug nombreux'
$ws.Range("D18").Value = 'Beta'
$ws.Range("E18").Value = 20130105
$ws.Range("F18").Value = 20140530
$ws.Range("H18").Value = 7
$ws.Range("A19").Value = 'MKart8ML02'
$ws.Range("B19").Value = 'Mario Kart 8 Multi Local 1.1'
$ws.Range("C19").Value = '2er version du Multi Local de Mario Kart 8, ajout mode course, mini-jeux, etc'
$ws.Range("D19").Value = 'Beta'
$ws.Range("E19").Value = 20140110
$ws.Range("F19").Value = 20140530
$ws.Range("H19").Value = 7
$ws.Range("A20").Value = 'MKart8MO01'
$ws.Range("B20").Value = 'Mario Kart 8 Multi online 1.0'
$ws.Range("C20").Value = '1er version du Multi en ligne de Mario Kart 8, bug nombreux'
$ws.Range("D20").Value = 'Beta'
$ws.Range("E20").Value = 20130105
$ws.Range("F20").Value = 20140530
$ws.Range("H20").Value = 8
$ws.Range("J20").Value = 'Max taille CodeVersion'
$ws.Range("A21").Value = 'MKart8MO02'
$ws.Range("B21").Value = 'Mario Kart 8 Multi online 1.1'
$ws.Range("C21").Value = '2er version du Multi en ligne de Mario Kart 8, ajout mode course, mini-jeux, etc'
$ws.Range("D21").Value = 'Beta'
$ws.Range("E21").Value = 20140110
$ws.Range("F21").Value = 20140530
$ws.Range("H21").Value = 8
$ws.Range("A22").Value = 'WOW01'
$ws.Range("B22").Value = 'world of warcraft 1.0'
$ws.Range("C22").Value = '1er version d''essaie du openworld, nombreux crash serveur et bug joueur'
$ws.Range("D22").Value = 'Alpha'
$ws.Range("E22").Value = 20001103
$ws.Range("F22").Value = 20050211
$ws.Range("H22").Value = 9
$ws.Range("A23").Value = 'WOW02'
$ws.Range("B23").Value = 'world of warcraft 1.1'
$ws.Range("C23").Value = '2e version, ajout des donjons, des évènement aléatoire, du craft et etc.'
$ws.Range("D23").Value = 'Beta'
$ws.Range("E23").Value = 20030115
$ws.Range("F23").Value = 20050211
$ws.Range("H23").Value = 9
$ws.Range("A24").Value = 'WOWCT01'
$ws.Range("B24").Value = 'world of warcraft CT  1.0'
$ws.Range("C24").Value = '1e version d''extension, Refonte du monde, nouvelle classe, race, lieux, etc.'
$ws.Range("D24").Value = 'Beta'
$ws.Range("E24").Value = 20081201
$ws.Range("F24").Value = 20101207
$ws.Range("H24").Value = 10
$ws.Range("A25").Value = 'NHL15SS01'
$ws.Range("B25").Value = 'NHL 15 solo 1.0'
$ws.Range("C25").Value = '1er version, nouveau joueur, peu de bug, ressemblance avec NHL 2014'
$ws.Range("D25").Value = 'Beta'
$ws.Range("E25").Value = 20130905
$ws.Range("F25").Value = 20140912
$ws.Range("H25").Value = 11
$ws.Range("A26").Value = 'NHL15SS02'
$ws.Range("B26").Value = 'NHL 15 solo 1.1'
$ws.Range("C26").Value = '2e version, nouveau chandaille, nouveau commentateur'
$ws.Range("D26").Value = 'Beta'
$ws.Range("E26").Value = 20131206
$ws.Range("F26").Value = 20140912
$ws.Range("H26").Value = 11
$ws.Range("A27").Value = 'NHL15ML01'
$ws.Range("B27").Value = 'NHL 15 Multi Local 1.0'
$ws.Range("C27").Value = '1er version, nouveau joueur, peu de bug, ressemblance avec NHL 2014'
$ws.Range("D27").Value = 'Beta'
$ws.Range("E27").Value = 20130905
$ws.Range("F27").Value = 20140912
$ws.Range("H27").Value = 12
$ws.Range("A28").Value = 'NHL15ML02'
$ws.Range("B28").Value = 'NHL 15 Multi Local 1.1'
$ws.Range("C28").Value = '2e version, nouveau chandaille, nouveau commentateur'
$ws.Range("D28").Value = 'Beta'
$ws.Range("E28").Value = 20131206
$ws.Range("F28").Value = 20140912
$ws.Range("H28").Value = 12
$ws.Range("A29").Value = 'NHL15MO01'
$ws.Range("B29").Value = 'NHL 15 Multi en ligne 1.0'
$ws.Range("C29").Value = '1er version, nouveau joueur, peu de bug, ressemblance avec NHL 2014'
$ws.Range("D29").Value = 'Beta'
$ws.Range("E29").Value = 20130905
$ws.Range("F29").Value = 20140912
$ws.Range("H29").Value = 13
$ws.Range("A30").Value = 'NHL15MO02'
$ws.Range("B30").Value = 'NHL 15 Multi en ligne 1.1'
$ws.Range("C30").Value = '2e version, nouveau chandaille, nouveau commentateur'
$ws.Range("D30").Value = 'Beta'
$ws.Range("E30").Value = 20131206
$ws.Range("F30").Value = 20140912
$ws.Range("H30").Value = 13
$ws.Range("A31").Value = 'NHL14SS01'
$ws.Range("B31").Value = 'NHL 14 solo 1.0'
$ws.Range("C31").Value = '1er version, nouveau joueur, peu de bug, ressemblance avec NHL 2013'
$ws.Range("D31").Value = 'Beta'
$ws.Range("E31").Value = 20120905
$ws.Range("F31").Value = 20130913
$ws.Range("H31").Value = 14
$ws.Range("A32").Value = 'NHL14SS02'
$ws.Range("B32").Value = 'NHL 14 solo 1.1'
$ws.Range("C32").Value = '2e version, nouveau chandaille, nouveau commentateur'
$ws.Range("D32").Value = 'Beta'
$ws.Range("E32").Value = 20121206
$ws.Range("F32").Value = 20130913
$ws.Range("H32").Value = 14
$ws.Range("A33").Value = 'NHL14ML01'
$ws.Range("B33").Value = 'NHL 14 Multi Local 1.0'
$ws.Range("C33").Value = '1er version, nouveau joueur, peu de bug, ressemblance avec NHL 2013'
$ws.Range("D33").Value = 'Beta'
$ws.Range("E33").Value = 20120905
$ws.Range("F33").Value = 20130913
$ws.Range("H33").Value = 15
$ws.Range("A34").Value = 'NHL14ML02'
$ws.Range("B34").Value = 'NHL 14 Multi Local 1.1'
$ws.Range("C34").Value = '2e version, nouveau chandaille, nouveau commentateur'
$ws.Range("D34").Value = 'Beta'
$ws.Range("E34").Value = 20121206
$ws.Range("F34").Value = 20130913
$ws.Range("H34").Value = 15
$ws.Range("A35").Value = 'NHL14MO01'
$ws.Range("B35").Value = 'NHL 14 Multi en ligne 1.0'
$ws.Range("C35").Value = '1er version, nouveau joueur, peu de bug, ressemblance avec NHL 2013'
$ws.Range("D35").Value = 'Beta'
$ws.Range("E35").Value = 20120905
$ws.Range("F35").Value = 20130913
$ws.Range("H35").Value = 16
$ws.Range("A36").Value = 'NHL14MO02'
$ws.Range("B36").Value = 'NHL 14 Multi en ligne 1.1'
$ws.Range("C36").Value = '2e version, nouveau chandaille, nouveau commentateur'
$ws.Range("D36").Value = 'Beta'
$ws.Range("E36").Value = 20121206
$ws.Range("F36").Value = 20130913
$ws.Range("H36").Value = 16
